$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new range of values used by the graph (G1:I1)
$ws.Range('G1').Value = 1
$ws.Range('H1').Value = 2
$ws.Range('I1').Value = 3

# Update C6 so it also depends on the new H1 value
$ws.Range('C6').Formula = '=A3+D2+H1'

# Register the defined name "test" pointing at the new range, for the graph links
$wb.Names.Add('test', '=Sheet1!$G$1:$I$1')
